$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get a numeric-looking new value (e.g. "212.72"). Excel would
# normally auto-convert such text to a real number on assignment, but the source
# data keeps these as literal text, so format them as Text first to preserve the
# string type, then restore the default (Normal) style afterwards.
$numericLikeCells = @("D5","D10","D15","D17","D20","D25","D37","D42","D43","D45","D47","D48")
foreach ($addr in $numericLikeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.233.58'
$ws.Range("D3").Value = '1.604.69'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '212.72'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").Value = '18.42'
$ws.Range("E10").Value = '  +2.25%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.827.39'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '1.607.09'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").Value = '0.512'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").Value = '26.205.79'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("D17").Value = '62.02'
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("E18").Value = '  +0.92%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '201.06'
$ws.Range("E20").Value = '  -2.16%  '
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  +2.70%  '
$ws.Range("D25").Value = '143.97'
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -1.90%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("E30").Value = '  +3.90%  '
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").Value = '1.161.41'
$ws.Range("E36").Value = '  +4.16%  '
$ws.Range("D37").Value = '0.0169'
$ws.Range("E37").Value = '  +3.64%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("E40").Value = '  +0.51%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.784'
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.32'
$ws.Range("E43").Value = '  +3.94%  '
$ws.Range("D44").Value = '1.739.58'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '91.71'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").Value = '0.0₆0106'
$ws.Range("E46").Value = '  +15.71%  '
$ws.Range("D47").Value = '1.53'
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("D48").Value = '54.13'
$ws.Range("E48").Value = '  +1.39%  '
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("E51").Value = '  -0.18%  '

# Restore default styling on the cells we temporarily reformatted as Text
foreach ($addr in $numericLikeCells) {
    $ws.Range($addr).Style = "Normal"
}
